# Update stats for 2026-01 (row 26: month 46023)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B26").Value = 6493
$ws.Range("D26").Value = 6053802
$ws.Range("E26").Value = 932.3582319420915
$ws.Range("F26").Value = 9.734662835896568
$ws.Range("H26").Value = 26.07063060076264
